$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns, placed right after the
# existing last column (AC). Set the text first, then clone the
# formatting (bold/bordered/centered header style) from the neighboring
# header cell AC1 via copy/paste-special so we reuse the same style index
# rather than creating a new one.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Team record (Wins/Losses/Ties) is the same for every player row since
# it describes the team, not the individual - 91-72-0 for every data row.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 91  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 72  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
